$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: TC_FB_011 - Password empty
$ws.Range("A13").Value = "TC_FB_011"
$ws.Range("B13").Value = "Password empty"
$ws.Range("C13").Value = "Navigate to facebook.com, Enter a valid format email, Enter a empty password click login"
$ws.Range("D13").Value = "Email valid@gmail.com and  password:empty"
$ws.Range("E13").Value = "Error msg :password is empty"
$ws.Range("F13").Value = "Error msg :password is empty"
$ws.Range("G13").Value = "pass"

# Row 14: TC_FB_012 - Email empty
$ws.Range("A14").Value = "TC_FB_012"
$ws.Range("B14").Value = "Email empty"
$ws.Range("C14").Value = "Navigate to facebook.com, Enter a ematy email, Enter a valid password click login"
$ws.Range("D14").Value = "Email empty@gmail.com and  password:valid"
$ws.Range("E14").Value = "Error msg :email is empty"
$ws.Range("F14").Value = "Error msg :email is empty"
$ws.Range("G14").Value = "pass"

# Match formatting of previous rows: wrap text style + row height 72
$ws.Range("C13:G14").WrapText = $true
$ws.Rows.Item(13).RowHeight = 72
$ws.Rows.Item(14).RowHeight = 72

# Update selection to match final state
$ws.Range("G14").Select()

# Touch page setup (portrait orientation), as occurs when the sheet is printed/previewed
$ws.PageSetup.Orientation = 1
